# Regenerate merged AHB files
# - rename the "_old" / "_new" suffixed headers to "_FV2210" / "_FV2304"
# - turn the sheet's used range into a proper Excel Table ("Table1")
# - freeze the header row (top row) so it stays visible while scrolling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells (row 1) from the old "_old"/"_new" suffix scheme to
#    the new "_FV2210"/"_FV2304" scheme.
$headerRenames = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# 2. Convert the used range A1:U56 into an Excel Table with headers.
$tableRange = $ws.Range("A1:U56")
$table = $ws.ListObjects.Add(1, $tableRange, $null, $true)
$table.Name = "Table1"

# 3. Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
